# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# de-de handback has completed and is now in sync with en-US:
#   - Status columns move from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - The "Latest Handback DateTime" timestamps advance
#   - The stale "Error Detail" messages (about an out-of-date handback file)
#     are cleared now that the handback is current
#   - A couple of report columns are widened so the (now longer) status /
#     datetime text is easier to read

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F)
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-22 16:50:41"
$wsZhCn.Range("P2").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-22 16:50:49"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
